$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35, shifting existing rows 35..56 down to 36..57.
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with a duplicate of the record that is
# now at row 36 (the original row 35 data), but with an updated Fecha (date)
# value, as described by the commit's weekly data update.
$ws.Cells.Item(35, 1).Value = 10
$ws.Cells.Item(35, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(35, 3).Value = "La Araucanía"
$ws.Cells.Item(35, 4).Value = 44663
$ws.Cells.Item(35, 5).Value = 9
$ws.Cells.Item(35, 6).Value = "Fruta"
$ws.Cells.Item(35, 7).Value = 100107
$ws.Cells.Item(35, 8).Value = "Otros"
$ws.Cells.Item(35, 9).Value = 100107011
$ws.Cells.Item(35, 10).Value = "Tuna"
$ws.Cells.Item(35, 11).Value = "Sin especificar"
$ws.Cells.Item(35, 12).Value = "Primera"
$ws.Cells.Item(35, 13).Value = 55
$ws.Cells.Item(35, 14).Value = 16000
$ws.Cells.Item(35, 15).Value = 16000
$ws.Cells.Item(35, 16).Value = 16000
$ws.Cells.Item(35, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(35, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(35, 19).Value = 1000
$ws.Cells.Item(35, 20).Value = 16
